# Update "想去人数" (want-to-go count) figures in column F on both the
# "展览" and "全部类型" worksheets, reflecting the newly generated data
# output (commit: "Update gh-pages to output generated at 456a3b4").

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

# Mapping of row number -> new value for column F
$updates = @{
    3  = 1731
    4  = 794
    7  = 12029
    10 = 481
    15 = 13517
    23 = 1913
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
